$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.535365
$ws.Range("H2").Value = 1.606095
$ws.Range("I2").Value = 0.1618182173563651
$ws.Range("J2").Value = 0.1618182173563651
$ws.Range("M2").Value = 1.028010333333333
$ws.Range("N2").Value = 3.084031
$ws.Range("O2").Value = 0.04165745457248912
$ws.Range("P2").Value = 0.04165745457248914
$ws.Range("Q2").Value = 0.5503607521049999
$ws.Range("R2").Value = 4.953246768945
$ws.Range("S2").Value = 0.00674093503852395
$ws.Range("T2").Value = 0.006740935038523951

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.535365
$ws.Range("H3").Value = 1.606095
$ws.Range("I3").Value = 0.1618182173563651
$ws.Range("J3").Value = 0.1618182173563651
$ws.Range("O3").Value = 0.9361395479363341
$ws.Range("P3").Value = 0.9361395479363344
$ws.Range("Q3").Value = 12.36788159442
$ws.Range("R3").Value = 111.31093434978
$ws.Range("S3").Value = 0.1514844328438511
$ws.Range("T3").Value = 0.1514844328438511

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.535365
$ws.Range("H4").Value = 1.606095
$ws.Range("I4").Value = 0.1618182173563651
$ws.Range("J4").Value = 0.1618182173563651
$ws.Range("O4").Value = 0.02220299749117665
$ws.Range("P4").Value = 0.02220299749117666
$ws.Range("Q4").Value = 0.293336655435
$ws.Range("R4").Value = 2.640029898915
$ws.Range("S4").Value = 0.003592849473990052
$ws.Range("T4").Value = 0.003592849473990053

$ws.Range("I5").Value = 0.6224306076670297
$ws.Range("J5").Value = 0.6224306076670296
$ws.Range("M5").Value = 1.028010333333333
$ws.Range("N5").Value = 3.084031
$ws.Range("O5").Value = 0.04165745457248912
$ws.Range("P5").Value = 0.04165745457248914
$ws.Range("Q5").Value = 2.116951867133667
$ws.Range("R5").Value = 19.052566804203
$ws.Range("S5").Value = 0.02592887476341609
$ws.Range("T5").Value = 0.0259288747634161

$ws.Range("I6").Value = 0.6224306076670297
$ws.Range("J6").Value = 0.6224306076670296
$ws.Range("O6").Value = 0.9361395479363341
$ws.Range("P6").Value = 0.9361395479363344
$ws.Range("S6").Value = 0.582681907683151
$ws.Range("T6").Value = 0.582681907683151

$ws.Range("I7").Value = 0.6224306076670297
$ws.Range("J7").Value = 0.6224306076670296
$ws.Range("O7").Value = 0.02220299749117665
$ws.Range("P7").Value = 0.02220299749117666
$ws.Range("S7").Value = 0.01381982522046262
$ws.Range("T7").Value = 0.01381982522046262

$ws.Range("I8").Value = 0.2157511749766052
$ws.Range("J8").Value = 0.2157511749766052
$ws.Range("M8").Value = 1.028010333333333
$ws.Range("N8").Value = 3.084031
$ws.Range("O8").Value = 0.04165745457248912
$ws.Range("P8").Value = 0.04165745457248914
$ws.Range("Q8").Value = 0.7337924052528887
$ws.Range("R8").Value = 6.604131647276
$ws.Range("S8").Value = 0.008987644770549082
$ws.Range("T8").Value = 0.008987644770549088

$ws.Range("I9").Value = 0.2157511749766052
$ws.Range("J9").Value = 0.2157511749766052
$ws.Range("O9").Value = 0.9361395479363341
$ws.Range("P9").Value = 0.9361395479363344
$ws.Range("S9").Value = 0.2019732074093321
$ws.Range("T9").Value = 0.2019732074093322

$ws.Range("I10").Value = 0.2157511749766052
$ws.Range("J10").Value = 0.2157511749766052
$ws.Range("O10").Value = 0.02220299749117665
$ws.Range("P10").Value = 0.02220299749117666
$ws.Range("S10").Value = 0.00479032279672398
$ws.Range("T10").Value = 0.004790322796723982
